# "Fixed bay update issues"
# - rows 2 & 3: longer Serial Number, Rank bumped 1 -> 2
# - row 6: Rank bumped 1 -> 2
# - three new bay rows appended (7, 8, 9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- existing row fixes --------------------------------------------------
$ws.Range("C2").Value = "adawddwdwadwad"
$ws.Range("E2").Value = 2

$ws.Range("C3").Value = "adawddwdwadwad"
$ws.Range("E3").Value = 2

$ws.Range("E6").Value = 2

# -- new row 7: bay 4-1 ----------------------------------------------------
$ws.Range("A7").Value = "4-1"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "esfe"
$ws.Range("D7").Value = "efsesf"
$ws.Range("E7").Value = 1
# F7 has no URL -- write it as an explicit empty text value (matches F6's
# existing empty-string cell) rather than leaving the cell untouched/blank.
$ws.Range("F7").Value = "'"
$ws.Range("F7").ClearFormats()

# -- new row 8: bay 6-2 ----------------------------------------------------
$ws.Range("A8").Value = "6-2"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "sefse"
$ws.Range("D8").Value = "efssef"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "'"
$ws.Range("F8").ClearFormats()

# -- new row 9: bay 6-1 ----------------------------------------------------
$ws.Range("A9").Value = "6-1"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "esef"
$ws.Range("D9").Value = "efsfes"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = "https://google.com"
